$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Champion")
$ws.Activate()

# --- New champion row 26: Mukai ---
# Copy the formatting of the previous data row (25) for column A and C:T,
# skipping column B (no "Skill" text for this champion yet) so no B26
# cell is emitted at all.
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C25:T25").Copy()
$ws.Range("C26:T26").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A26").Value = "Mukai"
$ws.Range("C26").Value = 100
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 10
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 1
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 1
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 2
$ws.Range("M26").Value = 1
$ws.Range("N26").Value = 0
$ws.Range("O26").Value = 3.8
$ws.Range("P26").Value = 0.25
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = 0
$ws.Range("S26").Value = 1.25
$ws.Range("T26").Value = 3

# --- New blank row 27 (A27 only, carries the same default cell style) ---
$ws.Range("A24").Copy()
$ws.Range("A27").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- View state: freeze-pane scroll + active selection, matching the author's edit ---
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$ws.Range("P27").Select()
